# SNAP.xlsx - "10Th - MB for single stock and added new group"
#
# The sheet tracks analyst ratings per broker (col A) with one column per
# "as of" date. This edit adds two new date columns (Jun_27, Jun_26) in
# front of the existing Jun_10 column, fills them with the placeholder
# rating "UN" for every existing broker row, and appends two new broker
# rows (Benchmark, Evercore ISI) for a newly tracked group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at B:C. This pushes the existing "Jun_10"
# data column from B to D (and its column-width formatting along with it).
$ws.Columns("B:C").Insert()

# New column headers.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"

# New column width for the inserted "Jun_26" column (C); column B keeps the
# default width, column D (old B) keeps its original bestFit width.
$ws.Columns("C:C").ColumnWidth = 10

# Fill the two new columns with the placeholder rating "UN" for every
# existing broker row (rows 2-27).
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Append the new broker group at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"

# Match the author's final selection.
$ws.Range("C7").Select()
